# The commit swaps the contents of ppt/theme/theme1.xml (the theme bound
# to the slide master / presentation, currently the "Integral" theme) and
# ppt/theme/theme2.xml (the theme bound only to the notes master,
# currently the default "Office Theme"). In this deck the font scheme
# (majorFont/minorFont) and the format scheme (fill/line/effect styles)
# are byte-identical between the two theme parts, so the only
# substantive difference is the 12-slot color scheme. Re-point the
# presentation's active theme (theme1.xml, reached via the slide
# master) at the "Office Theme" palette so the deck matches the target
# after the swap.

$p = $ppt.ActivePresentation
$master = $p.Slides.Item(1).Master
$colorScheme = $master.Theme.ThemeColorScheme

# RGB() style packed integers: R + G*256 + B*65536
$officeThemeColors = @{
    1  = 0         # dk1      000000
    2  = 16777215  # lt1      FFFFFF
    3  = 6968388   # dk2      44546A
    4  = 15132391  # lt2      E7E6E6
    5  = 13998939  # accent1  5B9BD5
    6  = 3243501   # accent2  ED7D31
    7  = 10855845  # accent3  A5A5A5
    8  = 49407     # accent4  FFC000
    9  = 12874308  # accent5  4472C4
    10 = 4697456   # accent6  70AD47
    11 = 12673797  # hlink    0563C1
    12 = 7491477   # folHlink 954F72
}

foreach ($idx in $officeThemeColors.Keys) {
    $colorScheme.Item($idx).RGB = $officeThemeColors[$idx]
}
